# Update the three-digit / one-digit division answers table.
# Each cell is addressed directly by (row, column) and its Range.Text is
# replaced in place so that formatting (rFonts/sz/jc) is preserved and so
# that Find/Execute cannot spill over into other cells that happen to share
# old/new text values (e.g. "913÷7=130, 3" is both an old and a new value).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
$cell.Range.Text = "338÷2=169, 0"
$cell = $t.Cell(1,2)
$cell.Range.Text = "549÷3=183, 0"
$cell = $t.Cell(1,3)
$cell.Range.Text = "768÷3=256, 0"
$cell = $t.Cell(1,4)
$cell.Range.Text = "565÷4=141, 1"
$cell = $t.Cell(1,5)
$cell.Range.Text = "211÷9=23, 4"
$cell = $t.Cell(5,1)
$cell.Range.Text = "309÷4=77, 1"
$cell = $t.Cell(5,2)
$cell.Range.Text = "850÷7=121, 3"
$cell = $t.Cell(5,3)
$cell.Range.Text = "801÷9=89, 0"
$cell = $t.Cell(5,4)
$cell.Range.Text = "352÷8=44, 0"
$cell = $t.Cell(5,5)
$cell.Range.Text = "968÷4=242, 0"
$cell = $t.Cell(9,1)
$cell.Range.Text = "492÷3=164, 0"
$cell = $t.Cell(9,2)
$cell.Range.Text = "412÷6=68, 4"
$cell = $t.Cell(9,3)
$cell.Range.Text = "230÷3=76, 2"
$cell = $t.Cell(9,4)
$cell.Range.Text = "913÷7=130, 3"
$cell = $t.Cell(9,5)
$cell.Range.Text = "796÷4=199, 0"
$cell = $t.Cell(13,1)
$cell.Range.Text = "811÷4=202, 3"
$cell = $t.Cell(13,2)
$cell.Range.Text = "320÷6=53, 2"
$cell = $t.Cell(13,3)
$cell.Range.Text = "493÷6=82, 1"
$cell = $t.Cell(13,4)
$cell.Range.Text = "514÷6=85, 4"
$cell = $t.Cell(13,5)
$cell.Range.Text = "658÷8=82, 2"
$cell = $t.Cell(17,1)
$cell.Range.Text = "691÷6=115, 1"
$cell = $t.Cell(17,2)
$cell.Range.Text = "675÷6=112, 3"
$cell = $t.Cell(17,3)
$cell.Range.Text = "234÷6=39, 0"
$cell = $t.Cell(17,4)
$cell.Range.Text = "860÷3=286, 2"
$cell = $t.Cell(17,5)
$cell.Range.Text = "194÷5=38, 4"
